$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the demo data to be more realistic: replace the placeholder
# "dataset_3" value with a real dataset name "dep_sante".
$ws.Range("B2").Value = "dep_sante"
$ws.Range("B3").Value = "dep_sante"

# Update the current selection to match the new state.
$ws.Range("B5").Select()
